$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6075
$ws.Range("E2").Value = 154
$ws.Range("F2").Value = 154
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 35
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 4281
$ws.Range("L2").Value = 2506
$ws.Range("M2").Value = 1774
$ws.Range("N2").Value = 1773
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 317
$ws.Range("Q2").Value = 231
$ws.Range("R2").Value = -34
$ws.Range("S2").Value = -171
$ws.Range("T2").Value = 63
$ws.Range("U2").Value = 168
$ws.Range("V2").Value = 1643
$ws.Range("W2").Value = 2.54
$ws.Range("X2").Value = 0.56
$ws.Range("Y2").Value = 1.97
$ws.Range("Z2").Value = 0.79
$ws.Range("AA2").Value = 141.25
$ws.Range("AB2").Value = 483.08
$ws.Range("AC2").Value = 55
$ws.Range("AD2").Value = 39.48
$ws.Range("AE2").Value = 3160
$ws.Range("AF2").Value = 0.68
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 63303370

# Row 3
$ws.Range("D3").Value = 6394
$ws.Range("E3").Value = 187
$ws.Range("F3").Value = 187
$ws.Range("G3").Value = 112
$ws.Range("H3").Value = 66
$ws.Range("I3").Value = 66
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4314
$ws.Range("L3").Value = 2449
$ws.Range("M3").Value = 1865
$ws.Range("N3").Value = 1864
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 317
$ws.Range("Q3").Value = 295
$ws.Range("R3").Value = -115
$ws.Range("S3").Value = -152
$ws.Range("T3").Value = 94
$ws.Range("U3").Value = 201
$ws.Range("V3").Value = 1530
$ws.Range("W3").Value = 2.93
$ws.Range("X3").Value = 1.04
$ws.Range("Y3").Value = 3.65
$ws.Range("Z3").Value = 1.55
$ws.Range("AA3").Value = 131.31
$ws.Range("AB3").Value = 511.81
$ws.Range("AC3").Value = 105
$ws.Range("AD3").Value = 16.93
$ws.Range("AE3").Value = 3321
$ws.Range("AF3").Value = 0.53
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 63303370

# Row 4
$ws.Range("D4").Value = 6401
$ws.Range("E4").Value = 150
$ws.Range("F4").Value = 150
$ws.Range("G4").Value = -35
$ws.Range("H4").Value = -58
$ws.Range("I4").Value = -58
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 4323
$ws.Range("L4").Value = 2513
$ws.Range("M4").Value = 1810
$ws.Range("N4").Value = 1808
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 317
$ws.Range("Q4").Value = 218
$ws.Range("R4").Value = -121
$ws.Range("S4").Value = -100
$ws.Range("T4").Value = 79
$ws.Range("U4").Value = 139
$ws.Range("V4").Value = 1495
$ws.Range("W4").Value = 2.35
$ws.Range("X4").Value = -0.9
$ws.Range("Y4").Value = -3.16
$ws.Range("Z4").Value = -1.34
$ws.Range("AA4").Value = 138.87
$ws.Range("AB4").Value = 492.95
$ws.Range("AC4").Value = -92
$ws.Range("AD4").Value = -17.43
$ws.Range("AE4").Value = 3223
$ws.Range("AF4").Value = 0.49
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 63303370

# Row 5
$ws.Range("D5").Value = 6399
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = -117
$ws.Range("H5").Value = -98
$ws.Range("I5").Value = -95
$ws.Range("J5").Value = -3
$ws.Range("K5").Value = 4216
$ws.Range("L5").Value = 2525
$ws.Range("M5").Value = 1691
$ws.Range("N5").Value = 1693
$ws.Range("O5").Value = -1
$ws.Range("P5").Value = 317
$ws.Range("Q5").Value = 13
$ws.Range("R5").Value = -56
$ws.Range("S5").Value = 63
$ws.Range("T5").Value = 87
$ws.Range("U5").Value = -74
$ws.Range("V5").Value = 1618
$ws.Range("W5").Value = 0.13
$ws.Range("X5").Value = -1.54
$ws.Range("Y5").Value = -5.44
$ws.Range("Z5").Value = -2.3
$ws.Range("AA5").Value = 149.26
$ws.Range("AB5").Value = 461.02
$ws.Range("AC5").Value = -150
$ws.Range("AD5").Value = -14.11
$ws.Range("AE5").Value = 3017
$ws.Range("AF5").Value = 0.7
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 63303370

# Row 6
$ws.Range("D6").Value = 6215
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 15
$ws.Range("G6").Value = -93
$ws.Range("H6").Value = -63
$ws.Range("I6").Value = -58
$ws.Range("K6").Value = 4123
$ws.Range("L6").Value = 2293
$ws.Range("M6").Value = 1829
$ws.Range("N6").Value = 1836
$ws.Range("P6").Value = 384
$ws.Range("Q6").Value = 50
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = -11
$ws.Range("T6").Value = 44
$ws.Range("U6").Value = 6
$ws.Range("V6").Value = 1445
$ws.Range("W6").Value = 0.24
$ws.Range("X6").Value = -1.02
$ws.Range("Y6").Value = -3.31
$ws.Range("Z6").Value = -1.52
$ws.Range("AA6").Value = 125.34
$ws.Range("AB6").Value = 402.85
$ws.Range("AC6").Value = -78
$ws.Range("AD6").Value = -23.88
$ws.Range("AE6").Value = 2640
$ws.Range("AF6").Value = 0.71
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 76738152

# Rows 7-9: clear all data columns (D:AJ), keep only A, B, C
$ws.Range("D7:AJ9").ClearContents()
